$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42845.5
$ws.Range("J3").Value = 42845.5
$ws.Range("L3").Value = 42845.5
$ws.Range("N3").Value = -43073.5
$ws.Range("H19").Value = 1668.0625
$ws.Range("I19").Value = 1298.7
$ws.Range("J19").Value = 2283.6667
$ws.Range("K19").Value = 1298.7
$ws.Range("L19").Value = 2283.6667
$ws.Range("M19").Value = -1123.7
$ws.Range("N19").Value = -2633.6667
$ws.Range("H28").Value = 1198.0526
$ws.Range("I28").Value = 1121.3529
$ws.Range("K28").Value = 1121.3529
$ws.Range("M28").Value = -636.3529000000001
$ws.Range("H32").Value = 1190.909
$ws.Range("J32").Value = 1199.8572
$ws.Range("L32").Value = 1199.8572
$ws.Range("N32").Value = -1851.8572
$ws.Range("H41").Value = 1033.7778
$ws.Range("I41").Value = 1207.4286
$ws.Range("K41").Value = 1207.4286
$ws.Range("M41").Value = -767.4286
$ws.Range("H98").Value = 5201.2
$ws.Range("I98").Value = 10000
$ws.Range("J98").Value = 4001.5
$ws.Range("K98").Value = 10000
$ws.Range("L98").Value = 4001.5
$ws.Range("M98").Value = -8502
$ws.Range("N98").Value = -6997.5
$ws.Range("H102").Value = 42845.5
$ws.Range("J102").Value = 42845.5
$ws.Range("L102").Value = 42845.5
$ws.Range("N102").Value = -49335.5
$ws.Range("H111").Value = 4000
$ws.Range("I111").Value = 5250
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 15750
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -12683
$ws.Range("N111").Value = -10634
$ws.Range("H122").Value = 5201.2
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 4001.5
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 12004.5
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -16904.5
$ws.Range("H135").Value = 21741580
$ws.Range("I135").Value = 2343.625
$ws.Range("J135").Value = 71431260
$ws.Range("K135").Value = 21092.625
$ws.Range("L135").Value = 642881340
$ws.Range("M135").Value = -18557.625
$ws.Range("N135").Value = -642886410
$ws.Range("H139").Value = 39140
$ws.Range("J139").Value = 39140
$ws.Range("L139").Value = 39140
$ws.Range("N139").Value = -49420
$ws.Range("H140").Value = 50500
$ws.Range("J140").Value = 50500
$ws.Range("L140").Value = 50500
$ws.Range("N140").Value = -60860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12860.585
$ws.Range("I32").Value = 13571
$ws.Range("J32").Value = 7745.6
$ws.Range("K32").Value = 13571
$ws.Range("L32").Value = 7745.6
$ws.Range("M32").Value = -13284
$ws.Range("N32").Value = -8319.6
$ws.Range("H45").Value = 2139.077
$ws.Range("I45").Value = 2040.8
$ws.Range("J45").Value = 2466.6667
$ws.Range("K45").Value = 2040.8
$ws.Range("L45").Value = 2466.6667
$ws.Range("M45").Value = -1663.8
$ws.Range("N45").Value = -3220.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4280.8057
$ws.Range("I105").Value = 3231.9
$ws.Range("J105").Value = 4684.231
$ws.Range("K105").Value = 3231.9
$ws.Range("L105").Value = 4684.231
$ws.Range("M105").Value = -1484.9
$ws.Range("N105").Value = -8178.231
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 34000
$ws.Range("J74").Value = 34000
$ws.Range("L74").Value = 34000
$ws.Range("N74").Value = -35748
$ws.Range("H77").Value = 34000
$ws.Range("J77").Value = 34000
$ws.Range("L77").Value = 102000
$ws.Range("N77").Value = -110736
$ws.Range("H96").Value = 39646.5
$ws.Range("J96").Value = 39646.5
$ws.Range("L96").Value = 39646.5
$ws.Range("N96").Value = -45138.5
$ws.Range("H122").Value = 1445.9667
$ws.Range("I122").Value = 1375.8182
$ws.Range("J122").Value = 1638.875
$ws.Range("K122").Value = 4127.4546
$ws.Range("L122").Value = 4916.625
$ws.Range("M122").Value = -1677.4546
$ws.Range("N122").Value = -9816.625
$ws.Range("H125").Value = 30158
$ws.Range("J125").Value = 30158
$ws.Range("L125").Value = 30158
$ws.Range("N125").Value = -35078
$ws.Range("H140").Value = 32100
$ws.Range("J140").Value = 32100
$ws.Range("L140").Value = 32100
$ws.Range("N140").Value = -42460

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 647.5
$ws.Range("I5").Value = 696
$ws.Range("J5").Value = 405
$ws.Range("K5").Value = 2088
$ws.Range("L5").Value = 1215
$ws.Range("M5").Value = -1976
$ws.Range("N5").Value = -1439
$ws.Range("H122").Value = 821
$ws.Range("I122").Value = 758.9583
$ws.Range("J122").Value = 1118.8
$ws.Range("K122").Value = 6830.6247
$ws.Range("L122").Value = 10069.2
$ws.Range("M122").Value = -4380.6247
$ws.Range("N122").Value = -14969.2
$ws.Range("H135").Value = 647.5
$ws.Range("I135").Value = 696
$ws.Range("J135").Value = 405
$ws.Range("K135").Value = 6264
$ws.Range("L135").Value = 3645
$ws.Range("M135").Value = -3729
$ws.Range("N135").Value = -8715

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2413.2
$ws.Range("I107").Value = 3150
$ws.Range("K107").Value = 3150
$ws.Range("M107").Value = -1230
$ws.Range("H132").Value = 6318.1577
$ws.Range("I132").Value = 5686.375
$ws.Range("J132").Value = 6777.636
$ws.Range("K132").Value = 17059.125
$ws.Range("L132").Value = 20332.908
$ws.Range("M132").Value = -14529.125
$ws.Range("N132").Value = -25392.908
$ws.Range("H138").Value = 56524.5
$ws.Range("J138").Value = 56524.5
$ws.Range("L138").Value = 56524.5
$ws.Range("N138").Value = -66804.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4694.654
$ws.Range("I40").Value = 6724.7
$ws.Range("J40").Value = 3425.875
$ws.Range("K40").Value = 6724.7
$ws.Range("L40").Value = 3425.875
$ws.Range("M40").Value = -6588.7
$ws.Range("N40").Value = -3697.875
$ws.Range("H136").Value = 35726436
$ws.Range("I136").Value = 62503136
$ws.Range("J136").Value = 24166
$ws.Range("K136").Value = 187509408
$ws.Range("L136").Value = 72498
$ws.Range("M136").Value = -187506858
$ws.Range("N136").Value = -77598
$ws.Range("H139").Value = 39910.445
$ws.Range("J139").Value = 39818
$ws.Range("L139").Value = 39818
$ws.Range("N139").Value = -50098

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H116").Value = 35680
$ws.Range("J116").Value = 35680
$ws.Range("L116").Value = 35680
$ws.Range("N116").Value = -44858
$ws.Range("H136").Value = 1267.138
$ws.Range("I136").Value = 1193.1305
$ws.Range("J136").Value = 1550.8334
$ws.Range("K136").Value = 3579.3915
$ws.Range("L136").Value = 4652.5002
$ws.Range("M136").Value = -1029.3915
$ws.Range("N136").Value = -9752.5002
